{"js": "// Replace the three-digit x one-digit multiplication expressions with\n// their updated values, as described by the diff.\nconst replacements = [\n  [\"343\u00d72=686\", \"544\u00d74=2176\"],\n  [\"860\u00d73=2580\", \"984\u00d76=5904\"],\n  [\"763\u00d73=2289\", \"553\u00d74=2212\"],\n  [\"683\u00d78=5464\", \"293\u00d79=2637\"],\n  [\"706\u00d75=3530\", \"861\u00d73=2583\"],\n  [\"632\u00d77=4424\", \"139\u00d74=556\"],\n  [\"586\u00d78=4688\", \"657\u00d78=5256\"],\n  [\"429\u00d79=3861\", \"307\u00d79=2763\"],\n  [\"729\u00d74=2916\", \"281\u00d75=1405\"],\n  [\"870\u00d72=1740\", \"845\u00d74=3380\"],\n  [\"296\u00d75=1480\", \"160\u00d74=640\"],\n  [\"628\u00d73=1884\", \"500\u00d78=4000\"],\n  [\"815\u00d76=4890\", \"446\u00d79=4014\"],\n  [\"685\u00d75=3425\", \"420\u00d74=1680\"],\n  [\"636\u00d78=5088\", \"900\u00d74=3600\"],\n  [\"236\u00d77=1652\", \"154\u00d74=616\"],\n  [\"102\u00d72=204\", \"877\u00d76=5262\"],\n  [\"770\u00d76=4620\", \"649\u00d76=3894\"],\n  [\"382\u00d73=1146\", \"796\u00d75=3980\"],\n  [\"945\u00d77=6615\", \"355\u00d75=1775\"],\n  [\"544\u00d75=2720\", \"203\u00d79=1827\"],\n  [\"585\u00d75=2925\", \"602\u00d79=5418\"],\n  [\"569\u00d75=2845\", \"122\u00d72=244\"],\n  [\"975\u00d76=5850\", \"241\u00d75=1205\"],\n  [\"394\u00d76=2364\", \"738\u00d79=6642\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the three-digit x one-digit multiplication expressions with\n# their updated values, as described by the diff.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @{old=\"343\u00d72=686\"; new=\"544\u00d74=2176\"},\n  @{old=\"860\u00d73=2580\"; new=\"984\u00d76=5904\"},\n  @{old=\"763\u00d73=2289\"; new=\"553\u00d74=2212\"},\n  @{old=\"683\u00d78=5464\"; new=\"293\u00d79=2637\"},\n  @{old=\"706\u00d75=3530\"; new=\"861\u00d73=2583\"},\n  @{old=\"632\u00d77=4424\"; new=\"139\u00d74=556\"},\n  @{old=\"586\u00d78=4688\"; new=\"657\u00d78=5256\"},\n  @{old=\"429\u00d79=3861\"; new=\"307\u00d79=2763\"},\n  @{old=\"729\u00d74=2916\"; new=\"281\u00d75=1405\"},\n  @{old=\"870\u00d72=1740\"; new=\"845\u00d74=3380\"},\n  @{old=\"296\u00d75=1480\"; new=\"160\u00d74=640\"},\n  @{old=\"628\u00d73=1884\"; new=\"500\u00d78=4000\"},\n  @{old=\"815\u00d76=4890\"; new=\"446\u00d79=4014\"},\n  @{old=\"685\u00d75=3425\"; new=\"420\u00d74=1680\"},\n  @{old=\"636\u00d78=5088\"; new=\"900\u00d74=3600\"},\n  @{old=\"236\u00d77=1652\"; new=\"154\u00d74=616\"},\n  @{old=\"102\u00d72=204\"; new=\"877\u00d76=5262\"},\n  @{old=\"770\u00d76=4620\"; new=\"649\u00d76=3894\"},\n  @{old=\"382\u00d73=1146\"; new=\"796\u00d75=3980\"},\n  @{old=\"945\u00d77=6615\"; new=\"355\u00d75=1775\"},\n  @{old=\"544\u00d75=2720\"; new=\"203\u00d79=1827\"},\n  @{old=\"585\u00d75=2925\"; new=\"602\u00d79=5418\"},\n  @{old=\"569\u00d75=2845\"; new=\"122\u00d72=244\"},\n  @{old=\"975\u00d76=5850\"; new=\"241\u00d75=1205\"},\n  @{old=\"394\u00d76=2364\"; new=\"738\u00d79=6642\"}\n)\n\nforeach ($p in $pairs) {\n  $find = $d.Content.Find\n  $find.Text = $p.old\n  $find.Replacement.Text = $p.new\n  $find.Execute($p.old, $false, $false, $false, $false, $false, $true, 1, $false, $p.new, 2) | Out-Null\n}\n"}
